$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.901.11'
$ws.Range("D3").Value = '1.564.87'
$ws.Range("E3").Value = '  +0.39%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '205.90'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.25%  '
$ws.Range("E6").Value = '  -1.02%  '
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '21.81'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.57%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.246'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.39%  '
$ws.Range("E10").Value = '  -1.03%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0864'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.23%  '
$ws.Range("D12").Value = '1.788.40'
$ws.Range("E12").Value = '  +0.46%  '
$ws.Range("D13").Value = '1.567.29'
$ws.Range("E13").Value = '  +0.48%  '
$ws.Range("E14").Value = '  -0.89%  '
$ws.Range("E15").Value = '  -0.17%  '
$ws.Range("D16").Value = '26.892.62'
$ws.Range("E16").Value = '  -0.92%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.31'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.48%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '215.44'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.15%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.36'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.88%  '
$ws.Range("D20").Value = '0.0₃0681'
$ws.Range("E20").Value = '  -0.58%  '
$ws.Range("E21").Value = '  -0.22%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.13'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.72%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.19'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.51%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.01'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.23%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.44'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.66%  '
$ws.Range("E26").Value = '  +2.17%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.95'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.30%  '
$ws.Range("E28").Value = '  -0.17%  '
$ws.Range("E29").Value = '  -0.74%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0465'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.90%  '
$ws.Range("E31").Value = '  -3.34%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.15'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.41%  '
$ws.Range("D33").Value = '1.394.95'
$ws.Range("E33").Value = '  +1.11%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.92'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.25%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.52'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.21%  '
$ws.Range("E36").Value = '  -0.52%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.922'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.30%  '
$ws.Range("E38").Value = '  -0.40%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.530'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.53%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.822'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.50%  '
$ws.Range("E41").Value = '  -0.15%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.993'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.17%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.52'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.80%  '
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("E45").Value = '  +1.25%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '63.80'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.86%  '
$ws.Range("D47").Value = '1.701.58'
$ws.Range("E47").Value = '  +0.55%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '86.75'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.57%  '
$ws.Range("D49").Value = '0.0₇0981'
$ws.Range("E49").Value = '  -0.14%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0501'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.07%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0954'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.58%  '
